# "Add files via upload" — Evan's row (row 8) on Sheet1 gets filled in with
# his skills / project-preference / e-commerce-experience / framework answers,
# the active selection moves to D8, and column C is widened to fit the new
# (longer) skills text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Evan's answers (row 8) — previously-blank cells under the
# "skills" / "would like to do on this project" / "Did you ever worked on
# E-Commerce website?" / "Framework knowledge" headers.
$ws.Range("C8").Value = "HTML(4), CSS(1-2), JavaScript(3-4), PHP(1-2)"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "N/A"
$ws.Range("D8").Value = "HTML or JavaScript"

# Column C needs to be widened so the longer skills text fits.
$ws.Columns.Item(3).ColumnWidth = 39.3

# Move the active selection to D8 (where the author was last working).
[void]$ws.Range("D8").Select()
